$wb = $excel.ActiveWorkbook

$rowData = @{
  1 = @("0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x01,0x18", "0xf", 380, "7.598631275147109e+23", 280, 15)
  2 = @("0x01,0x90", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x01,0x24", "0xe", 400, "5.68432987514711e+23", 292, 14)
  3 = @("0x00,0x6e", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x5F", "0x3", 110, "5.68631262647114e+23", 95, 3)
  4 = @("0x00,0x6e", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x5D", "0x3", 110, "9.85046333984776e+23", 93, 3)
}

for ($i = 1; $i -le 4; $i++) {
  $ws = $wb.Worksheets.Item($i)
  $newRow = 85
  $data = $rowData[$i]

  $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat
  $ws.Cells.Item($newRow, 1).Value = 45871.49369212963

  $ws.Cells.Item($newRow, 2).Value = $data[0]
  $ws.Cells.Item($newRow, 3).Value = $data[1]
  $ws.Cells.Item($newRow, 4).Value = $data[2]
  $ws.Cells.Item($newRow, 5).Value = $data[3]
  $ws.Cells.Item($newRow, 6).Value = $data[4]
  $ws.Cells.Item($newRow, 7).Value = [double]$data[5]
  $ws.Cells.Item($newRow, 8).Value = $data[6]
  $ws.Cells.Item($newRow, 9).Value = $data[7]
}
